# Update the rolling error-metric table with the newest ifoCAST evaluation
# quarter: drop the oldest row's values and shift the rest up by one,
# appending the freshly computed metrics for the newest quarter in the
# last row (labels in column A stay fixed; only B:G values move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.1143386297721907, 0.6889840097079422, 0.9692665404913205, 0.9845133521142924, 1.006200714633619, 18),
    @(-0.03231471279853988, 0.7695157101427289, 0.9504010021523444, 0.9748851225412891, 1.004336377456358, 17),
    @(0.05671622629984657, 0.6321910596607474, 0.8518300650585396, 0.9229464042177854, 0.951413452211165, 16),
    @(0.085928964333323, 0.5612671955956002, 0.5873845025919633, 0.7664101399328973, 0.7883078955770884, 15),
    @(0.1382138394267656, 0.7476276879241753, 1.097975515794379, 1.047843268716452, 1.077897322974905, 14),
    @(0.1215206328360093, 0.7819526689445239, 1.205421270947194, 1.097916786895616, 1.135726716235405, 13),
    @(0.2028461735751207, 0.8472735835465873, 1.158512388347609, 1.076342133500128, 1.10405820904526, 12),
    @(0.239907996146195, 0.6651045147971374, 0.5928596743010739, 0.769973814035954, 0.7673553351966808, 11),
    @(0.2225839694609408, 0.8306889768957216, 1.245987321698922, 1.116238021973325, 1.152988242321246, 10),
    @(0.2059614158583246, 1.002733503272531, 1.92968265145822, 1.389130178010045, 1.457110278930643, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
